$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4 / Column B (password) was mistakenly set to "iifEL1" instead of "iifEL".
# Fix it so the test data passes ("iifEL" matches rows 2 and 6).
$ws.Range("B4").Value = "iifEL"
